$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the Price (D) cells that are being updated as text, so
# numeric-looking strings (e.g. "314.94", "27.985.56") are not
# reinterpreted as numbers/dates by Excel. Only touch the cells
# that actually change so untouched cells keep their original style.
# (Cells are formatted one at a time -- a single comma-joined
# multi-area Range().NumberFormat assignment does not reliably
# apply to every area.)
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'

$ws.Range('D2').Value = '27.985.56'
$ws.Range('E2').Value = '  +3.32%  '
$ws.Range('D3').Value = '1.800.67'
$ws.Range('E3').Value = '  +3.89%  '
$ws.Range('D4').Value = '0.9978'
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').Value = '314.94'
$ws.Range('E5').Value = '  +1.31%  '
$ws.Range('D6').Value = '0.9976'
$ws.Range('E6').Value = '  -0.25%  '
$ws.Range('D7').Value = '0.5435'
$ws.Range('E7').Value = '  +11.33%  '
$ws.Range('D8').Value = '0.3787'
$ws.Range('E8').Value = '  +7.66%  '
$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').Value = '42.91'
$ws.Range('E9').Value = '  -0.46%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = '0.07565'
$ws.Range('E10').Value = '  +3.65%  '
$ws.Range('E11').Value = '  +7.32%  '
$ws.Range('D12').Value = '0.9970'
$ws.Range('E12').Value = '  -0.30%  '
$ws.Range('E13').Value = '  +5.32%  '
$ws.Range('D14').Value = '6.219'
$ws.Range('E14').Value = '  +5.39%  '
$ws.Range('D15').Value = '1.796.27'
$ws.Range('E15').Value = '  +3.76%  '
$ws.Range('D16').Value = '7.143'
$ws.Range('E16').Value = '  +3.36%  '
$ws.Range('D17').Value = '91.52'
$ws.Range('E17').Value = '  +4.62%  '
$ws.Range('D18').Value = '0.00001079'
$ws.Range('E18').Value = '  +3.69%  '
$ws.Range('D19').Value = '0.06507'
$ws.Range('E19').Value = '  +1.45%  '
$ws.Range('D20').Value = '0.9971'
$ws.Range('E20').Value = '  -0.28%  '
$ws.Range('E21').Value = '  +3.07%  '
$ws.Range('D22').Value = '5.966'
$ws.Range('E22').Value = '  +4.33%  '
$ws.Range('D23').Value = '28.006.23'
$ws.Range('E23').Value = '  +3.24%  '
$ws.Range('D24').Value = '11.23'
$ws.Range('E24').Value = '  +2.69%  '
$ws.Range('D25').Value = '2.093'
$ws.Range('E25').Value = '  +0.52%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '156.30'
$ws.Range('E26').Value = '  +1.21%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '20.59'
$ws.Range('E27').Value = '  +2.81%  '
$ws.Range('D28').Value = '2.402'
$ws.Range('E28').Value = '  +14.82%  '
$ws.Range('D29').Value = '2.003.55'
$ws.Range('E29').Value = '  +3.86%  '
$ws.Range('D30').Value = '122.53'
$ws.Range('E30').Value = '  +0.65%  '
$ws.Range('D31').Value = '1.146'
$ws.Range('E31').Value = '  +8.81%  '
$ws.Range('D32').Value = '0.1037'
$ws.Range('E32').Value = '  +11.10%  '
$ws.Range('D33').Value = '5.755'
$ws.Range('E33').Value = '  +6.03%  '
$ws.Range('D34').Value = '3.597'
$ws.Range('E34').Value = '  -1.37%  '
$ws.Range('D35').Value = '0.02298'
$ws.Range('E35').Value = '  +4.61%  '
$ws.Range('D36').Value = '8.689'
$ws.Range('E36').Value = '  +15.86%  '
$ws.Range('B37').Value = 'Algorand'
$ws.Range('C37').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D37').Value = '0.2109'
$ws.Range('E37').Value = '  +5.32%  '
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').Value = '5.036'
$ws.Range('E38').Value = '  +5.01%  '
$ws.Range('D39').Value = '0.06053'
$ws.Range('E39').Value = '  +1.77%  '
$ws.Range('E40').Value = '  +4.05%  '
$ws.Range('D41').Value = '0.6295'
$ws.Range('E41').Value = '  +4.55%  '
$ws.Range('D42').Value = '1.409'
$ws.Range('E42').Value = '  -1.59%  '
$ws.Range('D43').Value = '0.9969'
$ws.Range('E43').Value = '  -0.25%  '
$ws.Range('D44').Value = '1.152'
$ws.Range('E44').Value = '  +5.09%  '
$ws.Range('D45').Value = '13.37'
$ws.Range('E45').Value = '  +3.59%  '
$ws.Range('D46').Value = '0.5915'
$ws.Range('E46').Value = '  +3.97%  '
$ws.Range('D47').Value = '3.664'
$ws.Range('E47').Value = '  +2.14%  '
$ws.Range('D48').Value = '122.17'
$ws.Range('E48').Value = '  +2.74%  '
$ws.Range('D49').Value = '1.924'
$ws.Range('E49').Value = '  +3.77%  '
$ws.Range('D50').Value = '1.133'
$ws.Range('E50').Value = '  +2.27%  '
$ws.Range('E51').Value = '  +1.80%  '
